$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 5 (the "Resolving-Mac" target-cluster row) entirely.
$ws.Rows(5).Delete()

# Update recomputed TPM-derived values in rows 2-4, columns M:T.

# Row 2
$ws.Cells.Item(2, 13).Value2 = 0.377371
$ws.Cells.Item(2, 14).Value2 = 1.132113
$ws.Cells.Item(2, 15).Value2 = 0.4698794580655765
$ws.Cells.Item(2, 16).Value2 = 0.4698794580655764
$ws.Cells.Item(2, 17).Value2 = 0.03912230315066666
$ws.Cells.Item(2, 18).Value2 = 0.352100728356
$ws.Cells.Item(2, 19).Value2 = 0.4698794580655765
$ws.Cells.Item(2, 20).Value2 = 0.4698794580655764

# Row 3
$ws.Cells.Item(3, 15).Value2 = 0.443286188209444
$ws.Cells.Item(3, 16).Value2 = 0.443286188209444
$ws.Cells.Item(3, 18).Value2 = 0.3321732564800001
$ws.Cells.Item(3, 19).Value2 = 0.443286188209444
$ws.Cells.Item(3, 20).Value2 = 0.443286188209444

# Row 4
$ws.Cells.Item(4, 13).Value2 = 0.06973866666666667
$ws.Cells.Item(4, 14).Value2 = 0.209216
$ws.Cells.Item(4, 15).Value2 = 0.08683435372497944
$ws.Cells.Item(4, 16).Value2 = 0.08683435372497944
$ws.Cells.Item(4, 17).Value2 = 0.007229854065777779
$ws.Cells.Item(4, 18).Value2 = 0.06506868659200001
$ws.Cells.Item(4, 19).Value2 = 0.08683435372497944
$ws.Cells.Item(4, 20).Value2 = 0.08683435372497944
